$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.976.75"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.827.05"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'311.72"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.4642"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("D8").Value = "'0.3707"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "'0.07375"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'0.8732"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "'19.94"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "'0.07816"
$ws.Range("E12").Value = "  +6.30%  "
$ws.Range("D13").Value = "1.785.62"
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("D14").Value = "'6.591"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").Value = "'5.344"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "'91.67"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "'1.011"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "'0.000008838"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "'14.64"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "26.977.32"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("D22").Value = "'5.150"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").Value = "'10.58"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "2.023.24"
$ws.Range("E24").Value = "  -3.23%  "
$ws.Range("D25").Value = "'152.56"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'1.831"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D28").Value = "'2.081"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").Value = "'5.118"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").Value = "'115.62"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").Value = "'0.08861"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'2.974"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.444"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7269"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").Value = "'1.136"
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("D36").Value = "'2.493"
$ws.Range("E36").Value = "  +2.53%  "
$ws.Range("D37").Value = "'1.074"
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").Value = "'0.01950"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "'0.05221"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.928"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.212"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").Value = "'0.5196"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'0.8677"
$ws.Range("E43").Value = "  -13.99%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.1629"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").Value = "'8.216"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").Value = "'0.4836"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'10.21"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").Value = "'102.60"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").Value = "'1.625"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").Value = "'0.06215"
$ws.Range("E51").Value = "  -0.96%  "
